$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" (D) and "Volume(1h)" (E) columns store plain text values
# (e.g. "295.55", "1.78%") rather than numbers/percentages. Force each
# touched cell to Text format before writing so COM does not silently
# reinterpret the numeric-looking string as a Number/Percent value.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "295.55"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "1.78%"

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "31.07"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "0.83%"

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "4.907"
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "-0.83%"

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "0.07444"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "4.26%"

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "2.178"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "20.85%"

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "7.755"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "0.88%"

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "3.752"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "0.53%"

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.9136"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "1.88%"

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.08894"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "17.46%"

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.1711"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "3.95%"

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.08300"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "2.49%"

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.03158"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "3.15%"

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.1009"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "0.72%"

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.001523"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "1.32%"

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.005772"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "0.36%"

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "3.510"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "1.25%"

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.076"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "-0.20%"

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.3329"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "1.57%"

$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "-0.19%"

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "3.971"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "-1.67%"

$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "5.19%"

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.04556"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "0.75%"

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.001214"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "0.31%"

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.004625"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "15.60%"

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.0001300"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "4.19%"

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.0003400"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "-95.48%"

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01616"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "0.12%"

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.04479"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "2.65%"

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.007282"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "-0.60%"

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.008986"

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.1329"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "1.87%"

$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "-1.63%"

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.009144"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "0.29%"

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.00006123"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "1.78%"

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.00000000751"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "0.33%"

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.317"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "3.19%"

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.002003"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "-33.19%"

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.00002104"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "0.33%"

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0002004"
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "0.33%"
